$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Gloria Velez (row 2) left the list; Sebastian Romero (previously row 4)
# takes her place in row 2, and his old row is cleared out.
$ws.Range("A2").Value = "Sebastian Romero"
$ws.Range("B2").Value = 573138793438
$ws.Range("C2").Value = "Hno"

$ws.Range("A4:C4").ClearContents()

# Selection moved to C4 (now-empty row) in the saved view state.
$ws.Range("C4").Select() | Out-Null
